# Apply the "Updated cryptos list" data refresh (prices / 1h volume %)
# described by the commit diff. Column D holds numeric-looking strings
# (e.g. "209.51", "28.295.06") that must stay TEXT, exactly as authored -
# so each touched D cell is switched to text format ("@") before the
# value is written, then restored to the default "Normal" style so no
# stray formatting is introduced versus the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.295.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.550.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("E6").Value = '  -1.42%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.63'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.242'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.771.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.567.24'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.297.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.509'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("E20").Value = '  -2.57%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").Value = '  -3.14%  '
$ws.Range("E24").Value = '  -5.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.93%  '
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -3.18%  '
$ws.Range("E30").Value = '  -3.77%  '
$ws.Range("E31").Value = '  -4.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.381.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  +0.85%  '
$ws.Range("E36").Value = '  -3.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0161'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.509'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.89%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.90'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("E43").Value = '  -1.88%  '
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("E45").Value = '  -2.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("E47").Value = '  -6.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.685.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.22%  '
$ws.Range("E51").Value = '  +0.26%  '
